$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Tabelle4"
$newSheet.Range("B8").Value = "test"
$r = $newSheet.Range("B8")
$r.Interior.Pattern = -4105
Write-Host "set solid pattern only (no color)"
